$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet: Summary
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.55   # Current Capital
$summary.Range("B4").Value = 0.35      # Total P&L $
$summary.Range("B5").Value = 0.1       # Total P&L %
$summary.Range("B6").Value = 72        # Total Trades
$summary.Range("B8").Value = 29        # Losing Trades
$summary.Range("B9").Value = 45.83     # Win Rate %

# ------------------------------------------------------------------
# Sheet: Strategy Status (MarketMaking row, row 5)
# ------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.55     # Capital
$status.Range("D5").Value = 39         # Trades
$status.Range("E5").Value = 0.24       # P&L $
$status.Range("F5").Value = 0.55       # P&L %
$status.Range("G5").Value = 48.72      # Win Rate %

# ------------------------------------------------------------------
# Sheet: All Trades - close Trade #72 (row 73) and append Trade #105 (row 106)
# ------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G73").Value = 0.21
$allTrades.Range("H73").Value = "CLOSED"
$allTrades.Range("I73").Value = -12.5
$allTrades.Range("J73").Value = -0.03
$allTrades.Range("K73").Value = 100.55
$allTrades.Range("L73").Value = "early_exit"
$allTrades.Range("M73").Value = 0.11

# New row 106: clone the previous (structurally identical) row 105 so that
# text-like cells (dates) and blank cells retain the same cell types as the
# rest of the sheet, then overwrite the fields that actually differ.
$allTrades.Range("A105:Q105").Copy($allTrades.Range("A106:Q106"))
$allTrades.Range("A106").Value = 105
$allTrades.Range("C106").Value = "21:07:07"
$allTrades.Range("F106").Value = 0.24
$allTrades.Range("K106").Value = 100.5819219857093

# ------------------------------------------------------------------
# Sheet: MarketMaking - close Trade #72 (row 40) and append Trade #105 (row 73)
# ------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G40").Value = 0.21
$mm.Range("H40").Value = "CLOSED"
$mm.Range("I40").Value = -12.5
$mm.Range("J40").Value = -0.03
$mm.Range("K40").Value = 100.55
$mm.Range("P40").Value = "early_exit"
$mm.Range("Q40").Value = 0.11

# New row 73: clone the structurally identical row 41 (OPEN DOWN MarketMaking
# trade) so cell types for the date/blank columns match the sheet, then
# overwrite the differing fields.
$mm.Range("A41:Q41").Copy($mm.Range("A73:Q73"))
$mm.Range("A73").Value = 105
$mm.Range("C73").Value = "21:07:07"
$mm.Range("F73").Value = 0.24
$mm.Range("K73").Value = 100.5819219857093
